# Apply updates to column F (dSF) for specific rows as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    6  = 0
    7  = -2
    14 = 2
    24 = 0
    25 = 2
    27 = 0
    34 = 2
    41 = -1
    42 = 9
    43 = 4
    44 = -3
    52 = 1
    53 = -3
    57 = 2
    63 = -3
    65 = -1
    69 = -1
    73 = 0
    75 = 10
    85 = -1
    87 = 0
    90 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
